$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ID values to append after the existing 33 rows (header + 32 data rows)
$newIds = @("2430495", "2480050", "7914288", "9685871")

$startRow = 34
$endRow = $startRow + $newIds.Length - 1
$newRange = $ws.Range("A$startRow`:A$endRow")

# Force the values to be stored as text (matching the rest of the ID column,
# which keeps values like "0020893" intact instead of becoming numbers).
$newRange.NumberFormat = "@"
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $newIds[$i]
}

# Restore the exact cell style used by the rest of the data column (row 33).
$ws.Range("A33").Copy()
$newRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection to the newly added range, matching the diff
$newRange.Select()
